# "Changed the 'goto' statement on the EXCEL template"
#
# Sheet1 row 4 (the "What is the patient's sex?" question) has its
# Then_Goto (I4) / Else_Goto (J4) values switched from the text labels
# "pregnant" / "onset_date" to plain numeric row references 4 / 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5

# Leave the cursor where the author last left it when they saved.
$ws.Range("J6").Select()
